$wb = $excel.ActiveWorkbook

# Update the shared text "From Away" -> "Out of State" everywhere it is used
# (Region sheet category labels + the Region column on the main data sheet).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace('"From Away"', 'Out of State')
}

# Collapse the old selection on the "Region" sheet down to a single cell
# (B10) before leaving it.
$wsRegion = $wb.Worksheets.Item("Region")
$wsRegion.Range("B10").Select()

# Make "Portland Mayoral 2019" the active sheet/tab and leave the selection
# on L379.
$wsMain = $wb.Worksheets.Item("Portland Mayoral 2019")
$wsMain.Activate()
$wsMain.Range("L379").Select()
